# Fruta / hortaliza, semanal
# Insert the latest weekly record for "Ajo" (Vega Central Mapocho de Santiago)
# at row 332, pushing the existing historical rows (332-354) down to (333-355).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 332..354 down one row, freeing row 332 for the new record.
$ws.Rows.Item(332).Insert()

# Populate the new row 332 with the new weekly price observation.
$ws.Cells.Item(332, 1).Value  = 9
$ws.Cells.Item(332, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(332, 3).Value  = "Metropolitana"
$ws.Cells.Item(332, 4).Value  = 45166
$ws.Cells.Item(332, 5).Value  = 13
$ws.Cells.Item(332, 6).Value  = 100112003
$ws.Cells.Item(332, 7).Value  = "Ajo"
$ws.Cells.Item(332, 8).Value  = "Chino"
$ws.Cells.Item(332, 9).Value  = "Primera"
$ws.Cells.Item(332, 10).Value = 520
$ws.Cells.Item(332, 11).Value = 17000
$ws.Cells.Item(332, 12).Value = 18000
$ws.Cells.Item(332, 13).Value = 17500
$ws.Cells.Item(332, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(332, 15).Value = "China"
$ws.Cells.Item(332, 16).Value = 1750
$ws.Cells.Item(332, 17).Value = 10
$ws.Cells.Item(332, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells.Item(332, 4).NumberFormat = $ws.Cells.Item(333, 4).NumberFormat
